{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is the content of: async (context) => { ... }\n//\n// The edit rewrites the single summary paragraph:\n//  1. \"...in alpine areas with sharp microclimate gradients. Through a continuous\n//      germination experiment...\" becomes \"...in alpine areas with harsh climatic\n//      conditions and sharp fine-scale gradients. However, little is known about how\n//      these conditions influence germination timing. Through a continuous\n//      germination experiment...\"\n//  2. \"...two contrasting microclimatic conditions across a whole year. Using\n//      phenology traits we found a consistent germination phenological...\" becomes\n//      \"...two contrasting microclimatic conditions and found a consistent\n//      phenological...\"\n//  3. \"...shift. In warmer conditions, germination is anticipated between 45 and 60\n//      days with potential disrupting effects on plant communities.\" becomes\n//      \"...shift. Warmer conditions prompted earlier germination with potential\n//      disrupting effects on regeneration.\"\n\nconst body = context.document.body;\n\n// Each replacement below targets text that (per the document's original run\n// layout) lives inside a single run, so using Range.insertText(\u2026, \"Replace\")\n// rewrites that text in place without disturbing neighboring, unrelated runs.\nconst replacements = [\n  {\n    find: \" in alpine areas with sharp microclimate gradients. Through a continuous germination experiment\",\n    with: \" in alpine areas with harsh climatic conditions and sharp fine-scale gradients. However, little is known about how these conditions influence germination timing. Through a continuous germination experiment\",\n  },\n  {\n    find: \" two contrasting microclimatic conditions across a whole year. Using phenology traits we found a consistent germination phenological \",\n    with: \" two contrasting microclimatic conditions and found a consistent phenological \",\n  },\n  {\n    find: \". In warmer conditions, germination is anticipated between 45 and 60 days with potential disrupting effects on plant communities.\",\n    with: \". Warmer conditions prompted earlier germination with potential disrupting effects on regeneration.\",\n  },\n];\n\nfor (const { find, with: replacement } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + find);\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# The edit rewrites the single summary paragraph:\n#  1. \"...in alpine areas with sharp microclimate gradients. Through a continuous\n#      germination experiment...\" becomes \"...in alpine areas with harsh climatic\n#      conditions and sharp fine-scale gradients. However, little is known about how\n#      these conditions influence germination timing. Through a continuous\n#      germination experiment...\"\n#  2. \"...two contrasting microclimatic conditions across a whole year. Using\n#      phenology traits we found a consistent germination phenological...\" becomes\n#      \"...two contrasting microclimatic conditions and found a consistent\n#      phenological...\"\n#  3. \"...shift. In warmer conditions, germination is anticipated between 45 and 60\n#      days with potential disrupting effects on plant communities.\" becomes\n#      \"...shift. Warmer conditions prompted earlier germination with potential\n#      disrupting effects on regeneration.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $found = $find.Execute(\n        $findText,      # FindText\n        $false,         # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Replace-DocText: could not find expected text: $findText\"\n    }\n}\n\n# Apply the three replacements from the end of the paragraph toward the start.\n# (The underlying text-range rewrite re-serializes everything from the edit's\n# start position through to the end of the affected text each time it runs, so\n# working back-to-front keeps the not-yet-touched, earlier part of the\n# paragraph - \"Germination phenology ... survival, especially\" - intact.)\n\nReplace-DocText `\n    \". In warmer conditions, germination is anticipated between 45 and 60 days with potential disrupting effects on plant communities.\" `\n    \". Warmer conditions prompted earlier germination with potential disrupting effects on regeneration.\"\n\nReplace-DocText `\n    \" two contrasting microclimatic conditions across a whole year. Using phenology traits we found a consistent germination phenological \" `\n    \" two contrasting microclimatic conditions and found a consistent phenological \"\n\nReplace-DocText `\n    \" in alpine areas with sharp microclimate gradients. Through a continuous germination experiment\" `\n    \" in alpine areas with harsh climatic conditions and sharp fine-scale gradients. However, little is known about how these conditions influence germination timing. Through a continuous germination experiment\"\n\nWrite-Output $d.Content.Text\n"}
